$wb = $excel.ActiveWorkbook

# Add a new sheet "Out25" as a copy of "Set2025", placed after it (becomes sheetId=12, rId10)
$src = $wb.Worksheets.Item("Set2025")
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$src.Copy($null, $lastSheet)
$newSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet.Name = "Out25"

# Update the date in A1 to 2025-10-01 (Excel serial 45931)
$newSheet.Range("A1").Value = 45931

# Update the data block B2:M20 with the October figures
$data = New-Object 'object[,]' 19,12
$data[0,0] = 14713.62
$data[0,1] = 3597.2474999999999
$data[0,2] = 1235.1524999999999
$data[0,3] = 1235.1524999999999
$data[0,4] = 1235.1524999999999
$data[0,5] = 1852.72875
$data[0,6] = 1852.72875
$data[0,7] = 617.57624999999973
$data[0,8] = 617.57624999999973
$data[0,9] = 1235.1524999999995
$data[0,10] = 1235.1524999999995
$data[0,11] = 1003.59
$data[1,0] = 53765.829999999994
$data[1,1] = 6458.2937499999998
$data[1,2] = 4208.6076388888887
$data[1,3] = 4023.9476388888888
$data[1,4] = 4208.6076388888887
$data[1,5] = 8252.2983333333323
$data[1,6] = 8252.4383333333335
$data[1,7] = 3320.5661111111076
$data[1,8] = 3320.5661111111076
$data[1,9] = 5860.2522222222151
$data[1,10] = 5860.2522222222151
$data[1,11] = 3189.68
$data[2,0] = 5933.32
$data[2,1] = 741.66499999999996
$data[2,2] = 412.0361111111111
$data[2,3] = 412.0361111111111
$data[2,4] = 412.0361111111111
$data[2,5] = 988.88666666666654
$data[2,6] = 988.88666666666654
$data[2,7] = 329.62888888888847
$data[2,8] = 329.62888888888847
$data[2,9] = 659.25777777777694
$data[2,10] = 659.25777777777694
$data[2,11] = 488.25
$data[3,0] = 2748.8500000000004
$data[3,1] = 343.60625000000005
$data[3,2] = 190.89236111111111
$data[3,3] = 190.89236111111111
$data[3,4] = 190.89236111111111
$data[3,5] = 458.14166666666665
$data[3,6] = 458.14166666666665
$data[3,7] = 152.71388888888873
$data[3,8] = 152.71388888888873
$data[3,9] = 305.42777777777746
$data[3,10] = 305.42777777777746
$data[3,11] = 177.68
$data[4,0] = 2627.5
$data[4,1] = 328.4375
$data[4,2] = 182.46527777777777
$data[4,3] = 182.46527777777777
$data[4,4] = 182.46527777777777
$data[4,5] = 437.91666666666663
$data[4,6] = 437.91666666666663
$data[4,7] = 145.97222222222206
$data[4,8] = 145.97222222222206
$data[4,9] = 291.94444444444412
$data[4,10] = 291.94444444444412
$data[4,11] = 175.99
$data[5,0] = 2542.96
$data[5,1] = 317.87
$data[5,2] = 176.59444444444443
$data[5,3] = 176.59444444444443
$data[5,4] = 176.59444444444443
$data[5,5] = 423.82666666666665
$data[5,6] = 423.82666666666665
$data[5,7] = 141.2755555555554
$data[5,8] = 141.2755555555554
$data[5,9] = 282.5511111111108
$data[5,10] = 282.5511111111108
$data[5,11] = 157.18
$data[6,0] = 2043.06
$data[6,1] = 255.38249999999999
$data[6,2] = 141.87916666666666
$data[6,3] = 141.87916666666666
$data[6,4] = 141.87916666666666
$data[6,5] = 340.51
$data[6,6] = 340.51
$data[6,7] = 113.50333333333322
$data[6,8] = 113.50333333333322
$data[6,9] = 227.00666666666643
$data[6,10] = 227.00666666666643
$data[6,11] = 125.71
$data[7,0] = 2035.44
$data[7,1] = 254.43
$data[7,2] = 141.35
$data[7,3] = 141.35
$data[7,4] = 141.35
$data[7,5] = 339.23999999999995
$data[7,6] = 339.23999999999995
$data[7,7] = 113.07999999999988
$data[7,8] = 113.07999999999988
$data[7,9] = 226.15999999999977
$data[7,10] = 226.15999999999977
$data[7,11] = 133.33000000000001
$data[8,0] = 2858.98
$data[8,1] = 357.3725
$data[8,2] = 198.54027777777782
$data[8,3] = 198.54027777777782
$data[8,4] = 198.54027777777782
$data[8,5] = 476.49666666666661
$data[8,6] = 476.49666666666661
$data[8,7] = 158.83222222222207
$data[8,8] = 158.83222222222207
$data[8,9] = 317.66444444444414
$data[8,10] = 317.66444444444414
$data[8,11] = 253.78
$data[9,0] = -1863.76
$data[9,1] = -232.97
$data[9,2] = -129.42777777777778
$data[9,3] = -129.42777777777778
$data[9,4] = -129.42777777777778
$data[9,5] = -310.62666666666667
$data[9,6] = -310.62666666666667
$data[9,7] = -103.54222222222211
$data[9,8] = -103.54222222222211
$data[9,9] = -207.08444444444422
$data[9,10] = -207.08444444444422
$data[9,11] = 88.75
$data[10,0] = 0
$data[10,1] = 0
$data[10,2] = 0
$data[10,3] = 0
$data[10,4] = 0
$data[10,5] = 0
$data[10,6] = 0
$data[10,7] = 0
$data[10,8] = 0
$data[10,9] = 0
$data[10,10] = 0
$data[10,11] = 0
$data[11,0] = -524.93000000000006
$data[11,1] = -65.616250000000008
$data[11,2] = -36.453472222222224
$data[11,3] = -36.453472222222224
$data[11,4] = -36.453472222222224
$data[11,5] = -87.48833333333333
$data[11,6] = -87.48833333333333
$data[11,7] = -29.162777777777748
$data[11,8] = -29.162777777777748
$data[11,9] = -58.325555555555496
$data[11,10] = -58.325555555555496
$data[11,11] = 25
$data[12,0] = -527.49
$data[12,1] = -65.936250000000001
$data[12,2] = -36.631250000000001
$data[12,3] = -36.631250000000001
$data[12,4] = -36.631250000000001
$data[12,5] = -87.914999999999992
$data[12,6] = -87.914999999999992
$data[12,7] = -29.304999999999968
$data[12,8] = -29.304999999999968
$data[12,9] = -58.609999999999935
$data[12,10] = -58.609999999999935
$data[12,11] = 25.15
$data[13,0] = -2450.4199999999996
$data[13,1] = 0
$data[13,2] = -272.26888888888857
$data[13,3] = -272.26888888888857
$data[13,4] = -272.26888888888857
$data[13,5] = -408.4033333333341
$data[13,6] = -408.4033333333341
$data[13,7] = -136.13444444444428
$data[13,8] = -136.13444444444428
$data[13,9] = -272.26888888888857
$data[13,10] = -272.26888888888857
$data[13,11] = 116.7
$data[14,0] = 1147.3899999999981
$data[14,1] = 0
$data[14,2] = 127.48777777777744
$data[14,3] = 127.48777777777744
$data[14,4] = 127.48777777777744
$data[14,5] = 191.23166666666683
$data[14,6] = 191.23166666666683
$data[14,7] = 63.743888888888719
$data[14,8] = 63.743888888888719
$data[14,9] = 127.48777777777744
$data[14,10] = 127.48777777777744
$data[14,11] = 942.69
$data[15,0] = 15019.99
$data[15,1] = 0
$data[15,2] = 10013.326666666655
$data[15,3] = 0
$data[15,4] = 0
$data[15,5] = 2503.3316666666715
$data[15,6] = 2503.3316666666715
$data[15,7] = 0
$data[15,8] = 0
$data[15,9] = 0
$data[15,10] = 0
$data[15,11] = 928.77
$data[16,0] = -2669.37
$data[16,1] = -1334.6849999999999
$data[16,2] = -444.89499999999992
$data[16,3] = -444.89499999999992
$data[16,4] = -444.89499999999992
$data[16,5] = 0
$data[16,6] = 0
$data[16,7] = 0
$data[16,8] = 0
$data[16,9] = 0
$data[16,10] = 0
$data[16,11] = 0
$data[17,0] = 10504.48
$data[17,1] = 4980.8450000000003
$data[17,2] = 1841.211666666667
$data[17,3] = 1841.211666666667
$data[17,4] = 1841.211666666667
$data[17,5] = 0
$data[17,6] = 0
$data[17,7] = 0
$data[17,8] = 0
$data[17,9] = 0
$data[17,10] = 0
$data[17,11] = 723.64
$data[18,0] = 4584.54
$data[18,1] = 0
$data[18,2] = 2292.27
$data[18,3] = 2292.27
$data[18,4] = 0
$data[18,5] = 0
$data[18,6] = 0
$data[18,7] = 0
$data[18,8] = 0
$data[18,9] = 0
$data[18,10] = 0
$data[18,11] = 295.45999999999998

$newSheet.Range("B2:M20").Value = $data

# Set the active selection on the new sheet to A2 (matches the saved view state)
$newSheet.Range("A2").Select() | Out-Null
